$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns: h_x (m) and h_y (m), inserted right after "Elemento"
# (columns B and C were previously blank placeholder cells with style s=2).
$ws.Range("B1").Value = "h_x (m)"
$ws.Range("C1").Value = "h_y (m)"

# New data values for the two added columns (integer number format, like
# the "spt" column F which already uses a thousands-style integer format).
$ws.Range("B2").Value2 = 1
$ws.Range("C2").Value2 = 1
$ws.Range("B3").Value2 = 2
$ws.Range("C3").Value2 = 2

# Apply an integer ("#,##0") number format to the newly populated cells so
# they pick up a numFmtId=3 style (same numeric formatting family used
# elsewhere in the sheet), keeping the default "general" alignment.
$ws.Range("B2:C3").NumberFormat = "#,##0"
